$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking strings (e.g. "1.005", "26.940.24") that must
# remain literal text, matching the original inlineStr cells. Force text format
# before assignment so Excel does not auto-convert them to numbers.

$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = '26.940.24'
$ws.Cells.Item(2, 5).Value = '  -0.23%  '
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = '1.824.23'
$ws.Cells.Item(3, 5).Value = '  +0.09%  '
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = '1.005'
$ws.Cells.Item(4, 5).Value = '  -0.53%  '
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '310.56'
$ws.Cells.Item(5, 5).Value = '  +0.27%  '
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '1.003'
$ws.Cells.Item(6, 5).Value = '  -0.52%  '
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '0.4637'
$ws.Cells.Item(7, 5).Value = '  +0.14%  '
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '0.3690'
$ws.Cells.Item(8, 5).Value = '  +1.34%  '
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.07331'
$ws.Cells.Item(9, 5).Value = '  +0.50%  '
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '0.8750'
$ws.Cells.Item(10, 5).Value = '  +1.08%  '
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '0.07859'
$ws.Cells.Item(11, 5).Value = '  +3.35%  '
$ws.Cells.Item(12, 5).Value = '  -0.99%  '
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '1.861.78'
$ws.Cells.Item(13, 5).Value = '  -1.20%  '
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '5.331'
$ws.Cells.Item(14, 5).Value = '  +0.06%  '
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '6.546'
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '91.25'
$ws.Cells.Item(16, 5).Value = '  -2.28%  '
$ws.Cells.Item(17, 5).Value = '  -0.30%  '
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '0.000008816'
$ws.Cells.Item(18, 5).Value = '  +2.15%  '
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '1.002'
$ws.Cells.Item(19, 5).Value = '  -0.60%  '
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '14.75'
$ws.Cells.Item(20, 5).Value = '  +1.77%  '
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '26.962.18'
$ws.Cells.Item(21, 5).Value = '  -1.59%  '
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '5.097'
$ws.Cells.Item(22, 5).Value = '  -1.31%  '
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '10.52'
$ws.Cells.Item(23, 5).Value = '  -0.61%  '
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '2.022.80'
$ws.Cells.Item(24, 5).Value = '  -4.64%  '
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '151.86'
$ws.Cells.Item(25, 5).Value = '  +0.12%  '
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '1.858'
$ws.Cells.Item(26, 5).Value = '  +0.09%  '
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '18.40'
$ws.Cells.Item(27, 5).Value = '  +0.95%  '
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '2.034'
$ws.Cells.Item(28, 5).Value = '  -2.68%  '
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '5.098'
$ws.Cells.Item(29, 5).Value = '  -0.07%  '
$ws.Cells.Item(30, 5).Value = '  -0.49%  '
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '0.08871'
$ws.Cells.Item(31, 5).Value = '  -0.33%  '
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '2.958'
$ws.Cells.Item(32, 5).Value = '  +0.29%  '
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '0.7291'
$ws.Cells.Item(33, 5).Value = '  +0.07%  '
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '4.433'
$ws.Cells.Item(34, 5).Value = '  +0.11%  '
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '1.130'
$ws.Cells.Item(35, 5).Value = '  -1.00%  '
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '2.461'
$ws.Cells.Item(36, 5).Value = '  -1.00%  '
$ws.Cells.Item(37, 5).Value = '  -0.01%  '
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '0.01940'
$ws.Cells.Item(38, 5).Value = '  +1.10%  '
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '0.05212'
$ws.Cells.Item(39, 5).Value = '  -1.27%  '
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '2.955'
$ws.Cells.Item(40, 5).Value = '  +1.05%  '
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '7.077'
$ws.Cells.Item(41, 5).Value = '  -1.45%  '
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '0.5131'
$ws.Cells.Item(42, 5).Value = '  -1.45%  '
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '0.1623'
$ws.Cells.Item(43, 5).Value = '  -0.56%  '
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '8.141'
$ws.Cells.Item(44, 5).Value = '  -1.56%  '
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '0.4816'
$ws.Cells.Item(45, 5).Value = '  -1.07%  '
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '1.003'
$ws.Cells.Item(46, 5).Value = '  -0.55%  '
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '10.17'
$ws.Cells.Item(47, 5).Value = '  -0.04%  '
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '101.72'
$ws.Cells.Item(48, 5).Value = '  -1.52%  '
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '1.619'
$ws.Cells.Item(49, 5).Value = '  -0.82%  '
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '0.06201'
$ws.Cells.Item(50, 5).Value = '  -0.33%  '
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '64.51'
$ws.Cells.Item(51, 5).Value = '  -0.29%  '
